{"js": "// Update the \"three-digit number divided by one-digit number\" answer\n// table: each cell's division expression/answer text is replaced with\n// a newly generated one (same format: \"ABC\u00f7D=EF, G\").\nconst replacements = [\n  [\"658\u00f77=94, 0\", \"346\u00f73=115, 1\"],\n  [\"291\u00f74=72, 3\", \"891\u00f79=99, 0\"],\n  [\"615\u00f79=68, 3\", \"254\u00f76=42, 2\"],\n  [\"295\u00f78=36, 7\", \"980\u00f78=122, 4\"],\n  [\"849\u00f74=212, 1\", \"642\u00f73=214, 0\"],\n  [\"945\u00f77=135, 0\", \"208\u00f75=41, 3\"],\n  [\"565\u00f75=113, 0\", \"241\u00f73=80, 1\"],\n  [\"817\u00f77=116, 5\", \"611\u00f72=305, 1\"],\n  [\"485\u00f72=242, 1\", \"947\u00f73=315, 2\"],\n  [\"501\u00f78=62, 5\", \"979\u00f76=163, 1\"],\n  [\"943\u00f74=235, 3\", \"801\u00f74=200, 1\"],\n  [\"912\u00f76=152, 0\", \"746\u00f72=373, 0\"],\n  [\"605\u00f77=86, 3\", \"295\u00f77=42, 1\"],\n  [\"811\u00f77=115, 6\", \"916\u00f78=114, 4\"],\n  [\"652\u00f73=217, 1\", \"355\u00f79=39, 4\"],\n  [\"351\u00f77=50, 1\", \"748\u00f73=249, 1\"],\n  [\"312\u00f79=34, 6\", \"126\u00f74=31, 2\"],\n  [\"542\u00f76=90, 2\", \"336\u00f79=37, 3\"],\n  [\"650\u00f79=72, 2\", \"479\u00f75=95, 4\"],\n  [\"337\u00f72=168, 1\", \"847\u00f73=282, 1\"],\n  [\"647\u00f79=71, 8\", \"529\u00f72=264, 1\"],\n  [\"542\u00f72=271, 0\", \"554\u00f77=79, 1\"],\n  [\"463\u00f75=92, 3\", \"296\u00f77=42, 2\"],\n  [\"200\u00f75=40, 0\", \"668\u00f72=334, 0\"],\n  [\"522\u00f77=74, 4\", \"177\u00f79=19, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the \"three-digit number divided by one-digit number\" answer\n# table: each cell's division expression/answer text is replaced with\n# a newly generated one (same format: \"ABC\u00f7D=EF, G\").\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"658\u00f77=94, 0\", \"346\u00f73=115, 1\"),\n    @(\"291\u00f74=72, 3\", \"891\u00f79=99, 0\"),\n    @(\"615\u00f79=68, 3\", \"254\u00f76=42, 2\"),\n    @(\"295\u00f78=36, 7\", \"980\u00f78=122, 4\"),\n    @(\"849\u00f74=212, 1\", \"642\u00f73=214, 0\"),\n    @(\"945\u00f77=135, 0\", \"208\u00f75=41, 3\"),\n    @(\"565\u00f75=113, 0\", \"241\u00f73=80, 1\"),\n    @(\"817\u00f77=116, 5\", \"611\u00f72=305, 1\"),\n    @(\"485\u00f72=242, 1\", \"947\u00f73=315, 2\"),\n    @(\"501\u00f78=62, 5\", \"979\u00f76=163, 1\"),\n    @(\"943\u00f74=235, 3\", \"801\u00f74=200, 1\"),\n    @(\"912\u00f76=152, 0\", \"746\u00f72=373, 0\"),\n    @(\"605\u00f77=86, 3\", \"295\u00f77=42, 1\"),\n    @(\"811\u00f77=115, 6\", \"916\u00f78=114, 4\"),\n    @(\"652\u00f73=217, 1\", \"355\u00f79=39, 4\"),\n    @(\"351\u00f77=50, 1\", \"748\u00f73=249, 1\"),\n    @(\"312\u00f79=34, 6\", \"126\u00f74=31, 2\"),\n    @(\"542\u00f76=90, 2\", \"336\u00f79=37, 3\"),\n    @(\"650\u00f79=72, 2\", \"479\u00f75=95, 4\"),\n    @(\"337\u00f72=168, 1\", \"847\u00f73=282, 1\"),\n    @(\"647\u00f79=71, 8\", \"529\u00f72=264, 1\"),\n    @(\"542\u00f72=271, 0\", \"554\u00f77=79, 1\"),\n    @(\"463\u00f75=92, 3\", \"296\u00f77=42, 2\"),\n    @(\"200\u00f75=40, 0\", \"668\u00f72=334, 0\"),\n    @(\"522\u00f77=74, 4\", \"177\u00f79=19, 6\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
